# Re-pull data: update the dSF (F) column values to reflect the newly
# pulled/pushed data and recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    4  = 3
    5  = -5
    6  = -1
    7  = -1
    8  = 1
    9  = 10
    10 = 3
    11 = 2
    12 = -2
    13 = -4
    14 = 1
    15 = -3
    17 = -5
    18 = 3
    20 = 7
    21 = -4
    22 = -2
    23 = 7
    25 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
